$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-08-16 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-17 Thursday", 2) | Out-Null
$d.Content.Find.Execute("44÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "16÷3=", 2) | Out-Null
$d.Content.Find.Execute("71÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "56÷7=", 2) | Out-Null
$d.Content.Find.Execute("34÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "54÷2=", 2) | Out-Null
$d.Content.Find.Execute("37÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "55÷4=", 2) | Out-Null
$d.Content.Find.Execute("67÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "68÷6=", 2) | Out-Null
$d.Content.Find.Execute("38÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "47÷5=", 2) | Out-Null
$d.Content.Find.Execute("86÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "25÷4=", 2) | Out-Null
$d.Content.Find.Execute("79÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "22÷7=", 2) | Out-Null
$d.Content.Find.Execute("27÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "26÷9=", 2) | Out-Null
$d.Content.Find.Execute("76÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "68÷4=", 2) | Out-Null
$d.Content.Find.Execute("58÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "19÷4=", 2) | Out-Null
$d.Content.Find.Execute("49÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷8=", 2) | Out-Null
$d.Content.Find.Execute("36÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "93÷5=", 2) | Out-Null
$d.Content.Find.Execute("57÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "69÷6=", 2) | Out-Null
$d.Content.Find.Execute("77÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "78÷5=", 2) | Out-Null
$d.Content.Find.Execute("56÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "83÷6=", 2) | Out-Null
$d.Content.Find.Execute("26÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "77÷6=", 2) | Out-Null
$d.Content.Find.Execute("93÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "56÷5=", 2) | Out-Null
$d.Content.Find.Execute("13÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "58÷9=", 2) | Out-Null
$d.Content.Find.Execute("96÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "86÷4=", 2) | Out-Null
$d.Content.Find.Execute("72÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "50÷9=", 2) | Out-Null
$d.Content.Find.Execute("94÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "48÷7=", 2) | Out-Null
$d.Content.Find.Execute("24÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "22÷9=", 2) | Out-Null
$d.Content.Find.Execute("69÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷7=", 2) | Out-Null
$d.Content.Find.Execute("81÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "70÷6=", 2) | Out-Null
